# Adding 23/12/2017 exercise data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "WeightTraining"

# New workout rows for 23-Dec-2017 (Saturday), week 51, December 2017
# Columns: A=ExerciseId B=DateId C=Exercise Date D=Exercise Week E=Exercise Month
#          F=Exercise Year G=Exercise Day H=Exercise Name I=Weight J=Sets K=Reps
$rows = @(
    @(271, 33, 43092, 51, "December", 2017, "Saturday", "Dumbell Step Up",                    36,  4, 12),
    @(272, 33, 43092, 51, "December", 2017, "Saturday", "Bodyweight Dip",                     100, 4,  8),
    @(273, 33, 43092, 51, "December", 2017, "Saturday", "Bodyweight Pull-up",                 100, 5,  5),
    @(274, 33, 43092, 51, "December", 2017, "Saturday", "Squat Snatch",                         8, 4, 12),
    @(275, 33, 43092, 51, "December", 2017, "Saturday", "Box jumps",                            0, 3, 10),
    @(276, 33, 43092, 51, "December", 2017, "Saturday", "V-up crunches with medicine ball",     8, 4, 10),
    @(277, 33, 43092, 51, "December", 2017, "Saturday", "Leg Raises",                           0, 4, 10),
    @(278, 33, 43092, 51, "December", 2017, "Saturday", "Heel-taps",                            0, 4, 10),
    @(279, 33, 43092, 51, "December", 2017, "Saturday", "Barbell twists",                       0, 4, 10)
)

$startRow = 272
for ($n = 0; $n -lt $rows.Count; $n++) {
    $r = $startRow + $n
    $data = $rows[$n]
    $ws.Cells.Item($r, 1).Value  = $data[0]
    $ws.Cells.Item($r, 2).Value  = $data[1]
    $ws.Cells.Item($r, 3).Value  = $data[2]
    $ws.Cells.Item($r, 4).Value  = $data[3]
    $ws.Cells.Item($r, 5).Value  = $data[4]
    $ws.Cells.Item($r, 6).Value  = $data[5]
    $ws.Cells.Item($r, 7).Value  = $data[6]
    $ws.Cells.Item($r, 8).Value  = $data[7]
    $ws.Cells.Item($r, 9).Value  = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
    $ws.Cells.Item($r, 11).Value = $data[10]
}

# Column H (Exercise Name) needs to widen to fit the new longest entry
$ws.Columns.Item(8).ColumnWidth = 27.6667

# Scroll the frozen (bottom) pane down so the new rows are in view, then
# move the selection to the new last row, matching where the editor ended up
$win = $excel.ActiveWindow
$win.ScrollRow = 251
$win.ScrollColumn = 1
$ws.Range("B280").Select()
